# Estado de objetos ATM.xlsx - apply commit edits
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Selection moves to H20 ---
$null = $ws.Range("H20").Select()

# --- Create the new "black" fill/style by recoloring B19, then fan it out ---
$ws.Range("B19").Interior.ThemeColor = 1            # theme="1" (Black, Text 1) -> new fill + cellXf (s=8)

$null = $ws.Range("B19").Copy()
$null = $ws.Range("C19").PasteSpecial($xlPasteFormats)
$null = $ws.Range("B20").PasteSpecial($xlPasteFormats)
$null = $ws.Range("C20").PasteSpecial($xlPasteFormats)
$null = $ws.Range("B21").PasteSpecial($xlPasteFormats)
$null = $ws.Range("C21").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- Re-use the existing "green" style (s=4) on several cells ---
$null = $ws.Range("B6").Copy()
$null = $ws.Range("B8").PasteSpecial($xlPasteFormats)
$null = $ws.Range("B10").PasteSpecial($xlPasteFormats)
$null = $ws.Range("B11").PasteSpecial($xlPasteFormats)
$null = $ws.Range("B13").PasteSpecial($xlPasteFormats)
$null = $ws.Range("C13").PasteSpecial($xlPasteFormats)
$null = $ws.Range("B14").PasteSpecial($xlPasteFormats)
$null = $ws.Range("B15").PasteSpecial($xlPasteFormats)
$null = $ws.Range("B16").PasteSpecial($xlPasteFormats)
$null = $ws.Range("B17").PasteSpecial($xlPasteFormats)
$null = $ws.Range("B18").PasteSpecial($xlPasteFormats)
$null = $ws.Range("B23").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- Re-use the existing "yellow" style (s=5) ---
$null = $ws.Range("C3").Copy()
$null = $ws.Range("B7").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- Re-use the existing "orange" style (s=6) ---
$null = $ws.Range("B25").Copy()
$null = $ws.Range("C23").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- C9 gets B9's green style plus the "Listo" text ---
$null = $ws.Range("B9").Copy()
$null = $ws.Range("C9").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("C9").Value2 = "Listo"
